$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 182254.3736972502
$ws.Range("C2").Value = 151893.0518049554
$ws.Range("D2").Value = 62755.3576089744
$ws.Range("E2").Value = 114130.9399016776
$ws.Range("F2").Value = 34764.82301920529
$ws.Range("G2").Value = 69043.96513052752
$ws.Range("L2").Value = 119759.6573002258
$ws.Range("N2").Value = -2902623.09183835
$ws.Range("B3").Value = 342429.3398963325
$ws.Range("C3").Value = 1027641.30282144
$ws.Range("D3").Value = 219277.081236334
$ws.Range("E3").Value = 462164.0399398248
$ws.Range("F3").Value = 379963.4759659416
$ws.Range("G3").Value = 55064.36915698943
$ws.Range("L3").Value = 500077.3458097563
$ws.Range("N3").Value = -3202769.809684515
$ws.Range("B4").Value = 198376.4544156353
$ws.Range("C4").Value = 931364.2682337116
$ws.Range("D4").Value = 155738.5575611625
$ws.Range("E4").Value = 423272.1055467464
$ws.Range("F4").Value = 343793.9949011839
$ws.Range("G4").Value = 10500.1466368323
$ws.Range("L4").Value = 407205.8579777434
$ws.Range("N4").Value = -3649205.267098323
$ws.Range("B5").Value = 1153143.510197318
$ws.Range("C5").Value = 8583356.999522314
$ws.Range("D5").Value = 581369.2709010674
$ws.Range("E5").Value = 4473205.246548154
$ws.Range("F5").Value = 1110224.413824513
$ws.Range("G5").Value = 2036099.364222452
$ws.Range("L5").Value = 3555370.899690821
$ws.Range("N5").Value = -20899332.25343982
$ws.Range("B6").Value = 3073535.827925881
$ws.Range("C6").Value = 11986361.03160422
$ws.Range("D6").Value = 129592.9103568362
$ws.Range("E6").Value = 2005524.677776457
$ws.Range("F6").Value = 4072.4806395618
$ws.Range("G6").Value = 23224.77990687716
$ws.Range("H6").Value = 34881.47014419564
$ws.Range("I6").Value = 97575.03620374543
$ws.Range("L6").Value = 5030785.259659421
$ws.Range("N6").Value = -107069974.8612358
$ws.Range("B7").Value = 188897.8821044955
$ws.Range("C7").Value = 565496.7335720476
$ws.Range("D7").Value = 15103.38366958705
$ws.Range("E7").Value = 237786.7818006309
$ws.Range("H7").Value = 3438.630725903161
$ws.Range("I7").Value = 31512.84208492943
$ws.Range("L7").Value = 429080.3196151767
$ws.Range("N7").Value = -5512409.933681422
$ws.Range("B8").Value = 119773.535935784
$ws.Range("C8").Value = 836578.7910127238
$ws.Range("D8").Value = 13636.14454644389
$ws.Range("E8").Value = 323892.3258553279
$ws.Range("I8").Value = 38078.5820621863
$ws.Range("L8").Value = 554999.5459221107
$ws.Range("N8").Value = -3069857.358858988
$ws.Range("B9").Value = 202234.8876183375
$ws.Range("C9").Value = 1149890.666648876
$ws.Range("D9").Value = 30116.04890135345
$ws.Range("E9").Value = 407004.6630420301
$ws.Range("H9").Value = 6939.54417417106
$ws.Range("I9").Value = 51289.10659375732
$ws.Range("K9").Value = 102930.7367642276
$ws.Range("L9").Value = 772374.0044918191
$ws.Range("N9").Value = -2417439.500365156
